# Update leve-profit calculation cells (columns H-N) across multiple job sheets
# per scheduled market-data refresh. Values derived from updated market board pricing.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 92: Whinier than the Sword
$ws.Range("H92").Value = 1619.25
$ws.Range("I92").Value = 1619.25
$ws.Range("K92").Value = 1619.25
$ws.Range("M92").Value = -371.25
# Row 98: The Dotted Line
$ws.Range("H98").Value = 801.25
$ws.Range("I98").Value = 735
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 735
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 763
$ws.Range("N98").Value = -3996
# Row 107: Another Man's Ink
$ws.Range("H107").Value = 1345.3
$ws.Range("I107").Value = 1490.8889
$ws.Range("K107").Value = 1490.8889
$ws.Range("M107").Value = 429.1111000000001
# Row 122: Wishful Inking
$ws.Range("H122").Value = 801.25
$ws.Range("I122").Value = 735
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2205
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = 245
$ws.Range("N122").Value = -7900
# Row 132: Fast-forwarding Flora
$ws.Range("H132").Value = 1089.6923
$ws.Range("I132").Value = 1089.6923
$ws.Range("K132").Value = 3269.0769
$ws.Range("M132").Value = -739.0769
# Row 134: Binding Spells
$ws.Range("H134").Value = 150000
$ws.Range("J134").Value = 150000
$ws.Range("L134").Value = 150000
$ws.Range("N134").Value = -160140
# Row 135: For Tired Minds
$ws.Range("H135").Value = 2283.625
$ws.Range("I135").Value = 1086.6
$ws.Range("J135").Value = 4278.6665
$ws.Range("K135").Value = 9779.4
$ws.Range("L135").Value = 38507.9985
$ws.Range("M135").Value = -7244.4
$ws.Range("N135").Value = -43577.9985

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 3477.45
$ws.Range("I32").Value = 3449.9473
$ws.Range("K32").Value = 3449.9473
$ws.Range("M32").Value = -3162.9473
# Row 45: Hollow Hallmarks
$ws.Range("H45").Value = 2161.7273
$ws.Range("I45").Value = 2160.625
$ws.Range("K45").Value = 2160.625
$ws.Range("M45").Value = -1783.625
# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 4576.25
$ws.Range("I61").Value = 4691.5
$ws.Range("K61").Value = 4691.5
$ws.Range("M61").Value = -4479.5
# Row 97: Ore for Me
$ws.Range("H97").Value = 694.8333
$ws.Range("I97").Value = 561.1111
$ws.Range("K97").Value = 561.1111
$ws.Range("M97").Value = -65.11109999999996
# Row 122: Haste for High Durium
$ws.Range("H122").Value = 2589
$ws.Range("I122").Value = 2536.348
$ws.Range("K122").Value = 7609.044
$ws.Range("M122").Value = -5159.044
# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 4981.1665
$ws.Range("I132").Value = 4444.5
$ws.Range("J132").Value = 5249.5
$ws.Range("K132").Value = 13333.5
$ws.Range("L132").Value = 15748.5
$ws.Range("M132").Value = -10803.5
$ws.Range("N132").Value = -20808.5
# Row 136: Metal with Mettle
$ws.Range("H136").Value = 4576.25
$ws.Range("I136").Value = 4691.5
$ws.Range("K136").Value = 14074.5
$ws.Range("M136").Value = -11524.5

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 2760.2104
$ws.Range("I20").Value = 2780.3845
$ws.Range("J20").Value = 2716.5
$ws.Range("K20").Value = 2780.3845
$ws.Range("L20").Value = 2716.5
$ws.Range("M20").Value = -2533.3845
$ws.Range("N20").Value = -3210.5
# Row 86: Through Thick and Thin
$ws.Range("H86").Value = 2896.5
$ws.Range("I86").Value = 2844.75
$ws.Range("K86").Value = 2844.75
$ws.Range("M86").Value = -1721.75
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Range("H89").Value = 2896.5
$ws.Range("I89").Value = 2844.75
$ws.Range("K89").Value = 14223.75
$ws.Range("M89").Value = -8607.75
# Row 94: High Steal
$ws.Range("H94").Value = 987.2222
$ws.Range("I94").Value = 814.3333
$ws.Range("K94").Value = 814.3333
$ws.Range("M94").Value = -363.3333
# Row 99: Meddle in Metal
$ws.Range("H99").Value = 1000
$ws.Range("I99").Value = 1000
$ws.Range("K99").Value = 1000
$ws.Range("M99").Value = 498
# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 2785.2
$ws.Range("I105").Value = 2785.2
$ws.Range("K105").Value = 2785.2
$ws.Range("M105").Value = -1038.2
# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 1785.9131
$ws.Range("I134").Value = 1633
$ws.Range("J134").Value = 2512.25
$ws.Range("K134").Value = 4899
$ws.Range("L134").Value = 7536.75
$ws.Range("M134").Value = -2364
$ws.Range("N134").Value = -12606.75

$ws = $wb.Worksheets.Item("CRP")
# Row 16: Raise the Roof
$ws.Range("H16").Value = 6501.8
$ws.Range("I16").Value = 6992.6665
$ws.Range("J16").Value = 5765.5
$ws.Range("K16").Value = 6992.6665
$ws.Range("L16").Value = 5765.5
$ws.Range("M16").Value = -6705.6665
$ws.Range("N16").Value = -6339.5
# Row 42: Live Freelance or Die
$ws.Range("H42").Value = 6000
$ws.Range("I42").Value = 6000
$ws.Range("K42").Value = 6000
$ws.Range("M42").Value = -5407
# Row 107: Built to Last
$ws.Range("H107").Value = 1565.5
$ws.Range("I107").Value = 804.3333
$ws.Range("K107").Value = 804.3333
$ws.Range("M107").Value = 1115.6667
# Row 113: Patient Patients
$ws.Range("H113").Value = 6501.8
$ws.Range("I113").Value = 6992.6665
$ws.Range("J113").Value = 5765.5
$ws.Range("K113").Value = 6992.6665
$ws.Range("L113").Value = 5765.5
$ws.Range("M113").Value = -4822.6665
$ws.Range("N113").Value = -10105.5
# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 1494.8334
$ws.Range("I122").Value = 1494.8334
$ws.Range("K122").Value = 4484.5002
$ws.Range("M122").Value = -2034.5002
# Row 134: Wood You Be Quiet
$ws.Range("H134").Value = 3993.8333
$ws.Range("I134").Value = 4044.6
$ws.Range("J134").Value = 3740
$ws.Range("K134").Value = 12133.8
$ws.Range("L134").Value = 11220
$ws.Range("M134").Value = -9598.799999999999
$ws.Range("N134").Value = -16290

$ws = $wb.Worksheets.Item("GSM")
# Row 70: Sky Is the Limit
$ws.Range("H70").Value = 5419.6
$ws.Range("I70").Value = 5419.6
$ws.Range("K70").Value = 5419.6
$ws.Range("M70").Value = -5149.6
# Row 73: Hulls of Broken Dreams (L)
$ws.Range("H73").Value = 5419.6
$ws.Range("I73").Value = 5419.6
$ws.Range("K73").Value = 5419.6
$ws.Range("M73").Value = -4483.6
# Row 102: Put the Metal to the Peddle
$ws.Range("H102").Value = 3607.4167
$ws.Range("I102").Value = 3607.4167
$ws.Range("K102").Value = 3607.4167
$ws.Range("M102").Value = -1985.4167
# Row 113: Copious Crystal Cannons
$ws.Range("H113").Value = 2498.5
$ws.Range("I113").Value = 2498
$ws.Range("K113").Value = 2498
$ws.Range("M113").Value = -328
# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 998
$ws.Range("I122").Value = 998
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2994
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -544
$ws.Range("N122").ClearContents()
# Row 132: On Board for Lar
$ws.Range("H132").Value = 5997.8
$ws.Range("I132").Value = 5997.25
$ws.Range("K132").Value = 17991.75
$ws.Range("M132").Value = -15461.75

$ws = $wb.Worksheets.Item("WVR")
# Row 74: Clothing the Naked Truth
$ws.Range("H74").Value = 30000
$ws.Range("J74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31872
# Row 77: When in Robes (L)
$ws.Range("H77").Value = 30000
$ws.Range("J77").Value = 30000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -99360
# Row 107: Flax Wax
$ws.Range("H107").Value = 0
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws.Range("M107").ClearContents()
